{"js": "// Use case model added - update the \"software debe ser f\u00e1cil de utilizar\"\n// paragraph (\"intuitivo\" -> \"intuitivos\") and move the \"_GoBack\" bookmark\n// from the end of \"Asignaci\u00f3n de un plan alimenticio.\" up into this\n// paragraph (right after \"sin \").\n\n// 1) Remove the existing \"_GoBack\" bookmark (it currently sits right after\n//    \"Asignaci\u00f3n de un plan alimenticio.\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the target sentence and fix the wording: \"intuitivo\" -> \"intuitivos\".\nconst target = context.document.body.search(\n  \"personas sin conocimientos t\u00e9cnicos e intuitivo\",\n  { matchCase: true }\n);\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\nconst targetRange = target.items[0];\ntargetRange.insertText(\n  \"personas sin conocimientos t\u00e9cnicos e intuitivos\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark right after \"sin \" in the same sentence.\nconst anchor = context.document.body.search(\"personas sin \", { matchCase: true });\nanchor.load(\"text\");\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error(\"Bookmark anchor not found\");\n}\n\nconst anchorEnd = anchor.items[0].getRange(Word.RangeLocation.end);\nanchorEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Use case model added - update the \"software debe ser f\u00e1cil de utilizar\"\n# paragraph (\"intuitivo\" -> \"intuitivos\") and move the \"_GoBack\" bookmark\n# from the end of \"Asignaci\u00f3n de un plan alimenticio.\" up into this\n# paragraph (right after \"sin \").\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (currently right after\n#    \"Asignaci\u00f3n de un plan alimenticio.\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Fix the wording: \"intuitivo\" -> \"intuitivos\" in the target sentence.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"conocimientos t\u00e9cnicos e intuitivo;\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"conocimientos t\u00e9cnicos e intuitivos;\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 3) Re-insert the \"_GoBack\" bookmark right after \"sin \" in the same sentence.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"por personas sin \"\n$find2.Execute()\n$anchorRange = $find2.Parent\n$insertPoint = $d.Range($anchorRange.End, $anchorRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n"}
